$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H), mirroring the style used by the other
# header cells (B1:G1) in row 1.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Fill H2:H8 with the new "Save" values (all 0 for this sheet).
$ws.Range("H2:H8").Value = 0
